$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$formula = "=IF(Table23[[#This Row],[Given-When-Then (Tag)]]=`"`",IF(Table23[[#This Row],[Scenario]]=`"`",IF(Table23[[#This Row],[Feature]]<>`"`",CONCAT(`"Feature '`",Table23[[#This Row],[Feature]],`" `",Table23[[#This Row],[Sub Feature]],`"' {`"),`"`"),CONCAT(`"Scenario `",TEXT(Table23[[#This Row],[Scenario '#]],`"0000`"),`" '`",Table23[[#This Row],[Scenario]],`"' {`")),IF(INDIRECT(`"F`" & ROW() + 1)<>`"`",CONCAT(Table23[[#This Row],[Given-When-Then (Tag)]],`" '`",Table23[[#This Row],[Given-When-Then (Description)]],`"'`"),IF(INDIRECT(`"E`" & ROW() + 1)<>`"`",CONCAT(Table23[[#This Row],[Given-When-Then (Tag)]],`" '`",Table23[[#This Row],[Given-When-Then (Description)]],`"' }`"),CONCAT(Table23[[#This Row],[Given-When-Then (Tag)]],`" '`",Table23[[#This Row],[Given-When-Then (Description)]],`"' } }`"))))"
for ($r = 2; $r -le 17; $r++) {
    $ws2.Range("K$r").Formula = $formula
}
Write-Host "Done"
